# Adds the most recent shipment rows to the TAC tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Insert two new rows after the "BMOU2677749" row (old row 11, now row 11/12) ---
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

Set-RowValues 11 @("CCLU4722220", "EVER SHINE", "00207", "DJSEAA4007735", "7032014443", "6204914290")
Set-RowValues 12 @("TEMU6984979", "EVER SHINE", "01", "7032027243", "7032027243", "EGLV080900123743")

# --- Insert two new rows before the "TRLU6663841" row (now shifted to row 16) ---
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

Set-RowValues 16 @("SEGU5500923", "EVER SHINE", "01", "7032027246", "7032027246", "091930120641")
Set-RowValues 17 @("TEMU5580777", "EVER ENVOY", "00021", "DJSEAA4000159", "7032014429", "AOC0129056")

# --- Delete the obsolete "TCLU7869400" row (now at row 27) ---
$ws.Rows.Item(27).Delete()

# --- Insert a new row before the "FCIU7313352" row (now shifted to row 27) ---
$ws.Rows.Item(27).Insert()
Set-RowValues 27 @("CAIU6085579", "EVER SHINE", "00101", "DJSEAA4002110", "7032013727", "003900593568")

# --- Insert a new row before the "HMCU9058600" row (now shifted to row 36) ---
$ws.Rows.Item(36).Insert()
Set-RowValues 36 @("EGHU9313966", "EVER SHINE", "01", "7032027244", "7032027244", "EGLV149901010463")

# --- Append a new row at the very end (row 39) ---
Set-RowValues 39 @("FOLU3699332", "EVER SHINE", "00207", "DJSEAA4007741", "7032014443", "6204914290")

$wb.Save()
